$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 41, shifting rows 41..122 down to 42..123.
# This mirrors the edit in the source diff: a new weekly record was inserted
# before the existing "2022-12-28" record while keeping every other row intact
# (dimension grows from A1:R122 to A1:R123).
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new record's data.
$ws.Range("A41").Value2 = 10
$ws.Range("B41").Value = "Vega Modelo de Temuco"
$ws.Range("C41").Value = "La Araucanía"
$ws.Range("D41").Value2 = 45012
$ws.Range("E41").Value2 = 9
$ws.Range("F41").Value2 = 100112022
$ws.Range("G41").Value = "Arveja Verde"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value2 = 30
$ws.Range("K41").Value2 = 30000
$ws.Range("L41").Value2 = 30000
$ws.Range("M41").Value2 = 30000
$ws.Range("N41").Value = "$/saco 25 kilos"
$ws.Range("O41").Value = "Región de La Araucanía"
$ws.Range("P41").Value2 = 1200
$ws.Range("Q41").Value2 = 25
$ws.Range("R41").Value = "Hortaliza"
